$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (column D) and 1h volume change (column E) values,
# mirroring the GitHub Actions "Updated cryptos list" commit.
$data = @(
    @{ Row = 2; D = "23.244.05"; E = "  +0.98%  " },
    @{ Row = 3; D = "1.608.38"; E = "  +0.53%  " },
    @{ Row = 4; D = "1.000"; E = "  -0.12%  " },
    @{ Row = 5; D = "0.9995"; E = "  -0.17%  " },
    @{ Row = 6; D = "303.95"; E = "  +0.82%  " },
    @{ Row = 7; D = "0.3774"; E = "  -0.19%  " },
    @{ Row = 8; D = "52.07"; E = "  +4.73%  " },
    @{ Row = 9; D = "0.3640"; E = "  +0.03%  " },
    @{ Row = 10; D = "1.277"; E = "  +1.45%  " },
    @{ Row = 11; D = "0.08145"; E = "  +0.22%  " },
    @{ Row = 12; D = $null; E = "  -0.10%  " },
    @{ Row = 13; D = $null; E = "  +1.48%  " },
    @{ Row = 14; D = "6.597"; E = "  +0.12%  " },
    @{ Row = 15; D = $null; E = "  +0.93%  " },
    @{ Row = 16; D = "0.00001255"; E = "  +0.85%  " },
    @{ Row = 17; D = "1.606.40"; E = "  +0.07%  " },
    @{ Row = 18; D = "94.13"; E = "  +2.42%  " },
    @{ Row = 19; D = "0.06933"; E = "  +1.62%  " },
    @{ Row = 20; D = "18.18"; E = "  -0.15%  " },
    @{ Row = 21; D = "6.546"; E = "  +0.06%  " },
    @{ Row = 22; D = "1.000"; E = "  -0.13%  " },
    @{ Row = 23; D = $null; E = "  -1.80%  " },
    @{ Row = 24; D = "23.237.64"; E = "  +0.98%  " },
    @{ Row = 25; D = "2.446"; E = "  +3.93%  " },
    @{ Row = 26; D = "3.051"; E = "  +8.72%  " },
    @{ Row = 27; D = "21.24"; E = "  +0.90%  " },
    @{ Row = 28; D = "149.67"; E = "  -0.58%  " },
    @{ Row = 29; D = "5.286"; E = "  +1.06%  " },
    @{ Row = 30; D = "135.66"; E = "  +1.01%  " },
    @{ Row = 31; D = "2.381"; E = "  +2.93%  " },
    @{ Row = 32; D = "6.769"; E = "  -0.92%  " },
    @{ Row = 33; D = "1.775.68"; E = "  -0.54%  " },
    @{ Row = 34; D = "0.9676"; E = "  +0.52%  " },
    @{ Row = 35; D = "0.07507"; E = "  -1.06%  " },
    @{ Row = 36; D = "10.37"; E = "  +0.56%  " },
    @{ Row = 37; D = "0.02757"; E = "  +1.75%  " },
    @{ Row = 38; D = "0.2523"; E = "  -0.39%  " },
    @{ Row = 39; D = "6.142"; E = "  -1.84%  " },
    @{ Row = 40; D = "0.08812"; E = "  -0.91%  " },
    @{ Row = 41; D = "1.390"; E = "  +1.53%  " },
    @{ Row = 42; D = "0.7119"; E = "  +1.39%  " },
    @{ Row = 43; D = "12.50"; E = "  +0.68%  " },
    @{ Row = 44; D = "15.66"; E = "  +3.28%  " },
    @{ Row = 45; D = "0.6552"; E = "  -1.12%  " },
    @{ Row = 46; D = "2.335"; E = "  +1.39%  " },
    @{ Row = 47; D = "0.9984"; E = "  -0.19%  " },
    @{ Row = 48; D = "4.005"; E = "  +0.34%  " },
    @{ Row = 49; D = "133.00"; E = "  +1.15%  " },
    @{ Row = 50; D = "0.07955"; E = "  +0.61%  " },
    @{ Row = 51; D = $null; E = "  -1.84%  " }
)

foreach ($item in $data) {
    $r = $item.Row

    if ($item.D -ne $null) {
        $cellD = $ws.Cells.Item($r, 4)
        $styleD = $cellD.Style
        $cellD.NumberFormat = "@"
        $cellD.Value = $item.D
        $cellD.Style = $styleD
    }

    $cellE = $ws.Cells.Item($r, 5)
    $styleE = $cellE.Style
    $cellE.NumberFormat = "@"
    $cellE.Value = $item.E
    $cellE.Style = $styleE
}
